$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91: politeness_score (B91) was stored as text "3"; convert it to a
# genuine numeric value, matching the rest of the column.
$ws.Range("B91").Value = 3

# New row 92 appended at the end of the annotation table.
$ws.Range("A92").Value = "Ying Tang"
# B92 ("5") must stay a text value (matches the source data which mixes
# numeric-looking strings with real numbers) - leading apostrophe forces
# Excel to keep it as text instead of auto-converting to a number.
$ws.Range("B92").Value = "'5"
$ws.Range("C92").Value = " thank,thoughtful"
$ws.Range("D92").Value = "APC"
$ws.Range("E92").Value = "OTH"
$ws.Range("F92").Value = "c39fead7-b272-4988-9907-50ea12305918"
$ws.Range("G92").Value = "HknbyQbC-_annotated.xlsx"
$ws.Range("H92").Value = "We thank the reviewer for the thoughtful comments and suggestions."
